$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 64
$ws1.Range("F4").Value = 1481
$ws1.Range("F5").Value = 575
$ws1.Range("F6").Value = 1060
$ws1.Range("F7").Value = 10987
$ws1.Range("F8").Value = 10987
$ws1.Range("F13").Value = 749
$ws1.Range("F14").Value = 12199
$ws1.Range("F15").Value = 12714
$ws1.Range("F22").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 64
$ws4.Range("F5").Value = 1481
$ws4.Range("F6").Value = 575
$ws4.Range("F7").Value = 1060
$ws4.Range("F8").Value = 10987
$ws4.Range("F9").Value = 10987
$ws4.Range("F14").Value = 749
$ws4.Range("F15").Value = 12199
$ws4.Range("F16").Value = 12714
$ws4.Range("F23").Value = 26
